$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A44").Value = "평생학습자전형에 대해 알려줘."
$ws.Range("A45").Value = "평생학습자전형 지원자격은?"
$ws.Range("A46").Value = "실기고사는 어디로 가면 되나요?"
$ws.Range("A47").Value = "체육학과의 실기고사는 어디로 가면 되나요?"
$ws.Range("A49").Value = "2023학년도 경찰행정학과 수시모집 결과는?"
$ws.Range("A52").Value = "동의대학교 정시 일정에 대해 설명해주세요."
$ws.Range("A53").Value = "동의대학교 정시 모집 일정에 대해 설명해주세요."
$ws.Range("A58").Value = "농어촌학생전형으로 지원 시 제출해야할 서류는?"
$ws.Range("A59").Value = "농어촌학생전형으로 지원 시 12년 과정이 제출해야할 서류는?"
$ws.Range("A65").Value = "e비즈니스학전공 예비38번 이정도면 예비합격 가능할지 궁금합니다."
$ws.Range("A66").Value = "간호학과 - 동의대식 정시 성적 3.44인데 마지막으로 문 닫고 들어 올 가능성도 없나요??"
$ws.Range("A64").Value = "동의대학교 e비즈니스학과 정시 등급 5.3이면 최초합 가능성 있을까?"
$ws.Range("A67").Value = "간호학과 정시 성적 2.52 추가 합격으로도 합격 안되겠죠?"
$ws.Range("A69").Value = "간호학과 정시 성적이 2.58인데 정시 전형 합격할 수 있을까요?"
$ws.Range("A72").Value = "경영정보학과 정시 성적 3.75 인데 최초합으로 가능할까용?"
$ws.Range("A83").Value = "정시 성적으로 4.75 나오는데 가능할까요?? 그리고 유아교육학과는 이 성적으로 어렵겠죠 ?"
$ws.Range("A84").Value = "제가 동의대 정시 환산등급이 4.16인데 보육가정상담학과 일반학생전형으로 추합까지 생각해서 합격이 될까요?"
$ws.Range("A91").Value = "동의대 수능 정시 산출로 4.5나왔는데 추가모집이든 해서 들어갈수 있을까요?"
$ws.Range("A92").Value = "동의대학교 신소재 공학부 합격 가능할까요? 정시등급은 몇 인가여 신소재 학과는?"
$ws.Range("A93").Value = "제가 동의대학교 정시 표준점수가 447점으로 환산되는데 작년 신소재공학과 표점이 447점보다 1점 정도 위더라구요 (최종 평균)가능성은 있늘까요??"

$ws.Activate()
$ws.Range("A101").Select()
